$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# ------------------------------------------------------------------
# 1. Rename the existing "Walls+Windows+Lighting Package" (row 19) to
#    "Upgrade Package (Allow Individual Options)" - it keeps behaving
#    as the package that lets each option apply independently.
# ------------------------------------------------------------------
$ws.Range("B19").Value2 = "Upgrade Package (Allow Individual Options)"

# ------------------------------------------------------------------
# 2. Insert 6 new rows right before the "Building Characteristics
#    Report" measure block (old row 26) to make room for the new
#    "Upgrade Package (All or None)" package + its options + the new
#    package-level apply-logic argument.
# ------------------------------------------------------------------
$ws.Rows("26:31").Insert()

# ------------------------------------------------------------------
# 3. Row 26: new ApplyUpgrade measure header, mirrors row 19's layout.
# ------------------------------------------------------------------
$ws.Range("A19:X19").Copy()
$ws.Range("A26:X26").PasteSpecial(-4122)
$ws.Range("A26").Value2 = $true
$ws.Range("B26").Value2 = "Upgrade Package (All or None)"
$ws.Range("C26").Value2 = "ApplyUpgrade"
$ws.Range("D26").Value2 = "ApplyUpgrade"
$ws.Range("E26").Value2 = "RubyMeasure"

# ------------------------------------------------------------------
# 4. Row 27: new run_measure row, mirrors row 20's layout.
# ------------------------------------------------------------------
$cols27 = @("B","D","E","G","I","K","L","M","N","P","R")
foreach ($c in $cols27) {
    $ws.Range($c + "20").Copy()
    $ws.Range($c + "27").PasteSpecial(-4122)
}
$ws.Range("B27").Value2 = "variable"
$ws.Range("D27").Formula = '="Run " & B26'
$ws.Range("E27").Value2 = "run_measure"
$ws.Range("G27").Value2 = "Integer"
$ws.Range("I27").Value2 = 1
$ws.Range("K27").Value2 = 0
$ws.Range("L27").Value2 = 1
$ws.Range("M27").Value2 = 1
$ws.Range("N27").Value2 = 1
$ws.Range("P27").Value2 = "[0,1]"
$ws.Range("R27").Value2 = "discrete"

# ------------------------------------------------------------------
# 5. Rows 28-30: the three measure-argument rows (Option 1/2/3), same
#    shape as rows 21/23/25 (already inherit the correct row format
#    from the Insert above).
# ------------------------------------------------------------------
$ws.Range("A28").Value2 = ""
$ws.Range("B28").Value2 = "argument"
$ws.Range("D28").Value2 = "Option 1"
$ws.Range("E28").Formula = '=LOWER(SUBSTITUTE(D28," ","_"))'
$ws.Range("G28").Value2 = "string"
$ws.Range("I28").Value2 = "Insulation Wall|Wood Stud, R-13"

$ws.Range("B29").Value2 = "argument"
$ws.Range("D29").Value2 = "Option 2"
$ws.Range("E29").Formula = '=LOWER(SUBSTITUTE(D29," ","_"))'
$ws.Range("G29").Value2 = "string"
$ws.Range("I29").Value2 = "Windows|Low-E, Triple, Non-metal, Air, L-Gain"

$ws.Range("B30").Value2 = "argument"
$ws.Range("D30").Value2 = "Option 3"
$ws.Range("E30").Formula = '=LOWER(SUBSTITUTE(D30," ","_"))'
$ws.Range("G30").Value2 = "string"
$ws.Range("I30").Value2 = "Lighting|100% LED"

# ------------------------------------------------------------------
# 6. Row 31: new "Package Apply Logic" argument row.
# ------------------------------------------------------------------
$ws.Range("B31").Value2 = "argument"
$ws.Range("D31").Value2 = "Package Apply Logic"
$ws.Range("E31").Formula = '=LOWER(SUBSTITUTE(D31," ","_"))'
$ws.Range("G31").Value2 = "string"
$ws.Range("I31").Value2 = "Insulation Wall|Wood Stud, Uninsulated && (Windows|Clear, Single, Metal || Windows|Clear, Single, Non-metal)"

Write-Output "applied upgrade package apply-logic changes"
